# Update "想去人数" (number of people wanting to go) counts for several
# anime-convention events on both the "展览" sheet and the combined
# "全部类型" sheet, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value for column F
$exhibitUpdates = @{
    10 = 41
    11 = 7000
    12 = 252
    13 = 400
    14 = 3433
    15 = 247
    16 = 442
    19 = 58
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allUpdates = @{
    12 = 41
    14 = 7000
    16 = 252
    17 = 400
    18 = 3433
    19 = 247
    20 = 442
    23 = 58
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
